$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns D and E keep their original text representation
# (prices/percentages are stored as text, not numbers, in this sheet)
$ws.Columns("D:E").NumberFormat = "@"

$ws.Range("D2").Value = '26.620.54'
$ws.Range("D3").Value = '1.596.48'
$ws.Range("E3").Value = '  +0.47%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = '211.61'
$ws.Range("E5").Value = '  +0.28%  '
$ws.Range("E6").Value = '  +1.19%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("E8").Value = '  +0.38%  '
$ws.Range("D9").Value = '0.246'
$ws.Range("E9").Value = '  -0.21%  '
$ws.Range("D10").Value = '19.47'
$ws.Range("E10").Value = '  -0.52%  '
$ws.Range("E11").Value = '  +0.47%  '
$ws.Range("D12").Value = '1.820.83'
$ws.Range("E12").Value = '  +0.51%  '
$ws.Range("D13").Value = '1.573.17'
$ws.Range("E13").Value = '  -1.58%  '
$ws.Range("E14").Value = '  +0.22%  '
$ws.Range("E15").Value = '  +0.04%  '
$ws.Range("D16").Value = '64.59'
$ws.Range("E16").Value = '  -0.29%  '
$ws.Range("D17").Value = '26.606.50'
$ws.Range("E17").Value = '  -0.03%  '
$ws.Range("E18").Value = '  +0.60%  '
$ws.Range("D19").Value = '208.91'
$ws.Range("E19").Value = '  +0.45%  '
$ws.Range("E20").Value = '  +0.07%  '
$ws.Range("D21").Value = '6.96'
$ws.Range("E21").Value = '  +3.44%  '
$ws.Range("E22").Value = '  +0.50%  '
$ws.Range("D23").Value = '2.31'
$ws.Range("E23").Value = '  -1.21%  '
$ws.Range("D24").Value = '8.90'
$ws.Range("E24").Value = '  +0.55%  '
$ws.Range("D25").Value = '145.67'
$ws.Range("E25").Value = '  -0.74%  '
$ws.Range("E26").Value = '  -0.05%  '
$ws.Range("D27").Value = '7.11'
$ws.Range("E27").Value = '  -1.86%  '
$ws.Range("E28").Value = '  +0.77%  '
$ws.Range("D29").Value = '15.26'
$ws.Range("E29").Value = '  -0.06%  '
$ws.Range("D30").Value = '0.0508'
$ws.Range("E30").Value = '  +0.20%  '
$ws.Range("E31").Value = '  +0.49%  '
$ws.Range("E32").Value = '  +0.34%  '
$ws.Range("D33").Value = '0.655'
$ws.Range("E33").Value = '  -1.85%  '
$ws.Range("D34").Value = '2.94'
$ws.Range("E34").Value = '  +0.98%  '
$ws.Range("D35").Value = '1.283.86'
$ws.Range("E35").Value = '  -1.58%  '
$ws.Range("E36").Value = '  +1.12%  '
$ws.Range("E37").Value = '  +0.90%  '
$ws.Range("E38").Value = '  -0.39%  '
$ws.Range("D39").Value = '0.842'
$ws.Range("E39").Value = '  +1.68%  '
$ws.Range("E40").Value = '  +0.04%  '
$ws.Range("D41").Value = '5.49'
$ws.Range("E41").Value = '  +2.31%  '
$ws.Range("B42").Value = 'MXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D42").Value = '2.20'
$ws.Range("E42").Value = '  +1.80%  '
$ws.Range("B43").Value = 'Aave'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D43").Value = '64.45'
$ws.Range("E43").Value = '  +2.89%  '
$ws.Range("B44").Value = 'TrustWalletToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D44").Value = '0.785'
$ws.Range("E44").Value = '  -0.94%  '
$ws.Range("D45").Value = '1.733.23'
$ws.Range("E45").Value = '  +0.48%  '
$ws.Range("E46").Value = '  +8.96%  '
$ws.Range("D47").Value = '89.72'
$ws.Range("E47").Value = '  +0.18%  '
$ws.Range("D48").Value = '1.60'
$ws.Range("E48").Value = '  -0.61%  '
$ws.Range("E49").Value = '  +4.75%  '
$ws.Range("D50").Value = '0.0507'
$ws.Range("E50").Value = '  +0.48%  '
$ws.Range("D51").Value = '7.47'
$ws.Range("E51").Value = '  -0.18%  '
